$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: DATE_TYPE_CODE 002 -> 001
# Force text storage so the leading zero is preserved (not read as number 1),
# then clear the temporary formatting so no new cell style is left behind.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# N2: REPORT_DATE 2018-06-30 00:00:00 -> 2017-12-31 00:00:00
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# O2: TOTAL_ASSETS
$ws.Range("O2").Value = 2656390159.09

# P2: FIXED_ASSET
$ws.Range("P2").Value = 21354361.97

# Q2: MONETARYFUNDS
$ws.Range("Q2").Value = 7653575.7

# R2: MONETARYFUNDS_RATIO
$ws.Range("R2").Value = -92.9895729573

# S2: ACCOUNTS_RECE
$ws.Range("S2").Value = 977629111.45

# T2: ACCOUNTS_RECE_RATIO
$ws.Range("T2").Value = 52.5329209611

# U2: INVENTORY
$ws.Range("U2").Value = 1071828733.21

# V2: INVENTORY_RATIO
$ws.Range("V2").Value = 11.5628659955

# W2: TOTAL_LIABILITIES
$ws.Range("W2").Value = 1800616844.57

# X2: ACCOUNTS_PAYABLE
$ws.Range("X2").Value = 538089406.1799999

# Y2: ACCOUNTS_PAYABLE_RATIO
$ws.Range("Y2").Value = -22.6445008783

# Z2: ADVANCE_RECEIVABLES
$ws.Range("Z2").Value = 3930290.4

# AA2: ADVANCE_RECEIVABLES_RATIO
$ws.Range("AA2").Value = -38.7541736356

# AB2: TOTAL_EQUITY
$ws.Range("AB2").Value = 855773314.52

# AC2: TOTAL_EQUITY_RATIO
$ws.Range("AC2").Value = 38.3327316211

# AD2: TOTAL_ASSETS_RATIO
$ws.Range("AD2").Value = 29.9212022715

# AE2: TOTAL_LIAB_RATIO
$ws.Range("AE2").Value = 26.2720280472

# AF2: CURRENT_RATIO
$ws.Range("AF2").Value = 136.808642016

# AG2: DEBT_ASSET_RATIO
$ws.Range("AG2").Value = 67.7843515723
